$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 47 (A1:R47). New weekly data for
# row 46/47 has arrived; the previous values for those two rows are
# pushed down to new rows 48/49, and rows 46/47 are overwritten with
# the new figures.

# 1) Duplicate the existing rows 46:47 down to rows 48:49 before
#    overwriting them with the new data.
$src = $ws.Range("A46:R47")
$dst = $ws.Range("A48:R49")
$src.Copy($dst)

# 2) Overwrite row 46 with the new weekly figures.
$ws.Range("D46").Value = 44706
$ws.Range("J46").Value = 200
$ws.Range("K46").Value = 13000
$ws.Range("L46").Value = 14000
$ws.Range("M46").Value = 13500
$ws.Range("P46").Value = 750

# 3) Overwrite row 47 with the new weekly figures.
$ws.Range("D47").Value = 44706
$ws.Range("J47").Value = 100
